# The commit removes the explicit "contextual spacing" paragraph-format
# override (<w:contextualSpacing w:val="0"/>) that had been stamped onto
# every paragraph's pPr (in the main story and inside every comment) by
# the previous conversion pass. The correct, idiomatic Word automation
# call for that is ParagraphFormat.ContextualSpacing = False on each
# paragraph's range, so walk every paragraph in the document body and in
# every comment and clear it there.

$d = $word.ActiveDocument

foreach ($para in $d.Paragraphs) {
    $para.Range.ParagraphFormat.ContextualSpacing = $false
}

$commentCount = $d.Comments.Count
for ($i = 1; $i -le $commentCount; $i++) {
    $comment = $d.Comments.Item($i)
    foreach ($para in $comment.Range.Paragraphs) {
        $para.Range.ParagraphFormat.ContextualSpacing = $false
    }
}
